$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 109, shifting existing rows 109-113 down to 110-114.
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new weekly price record.
$ws.Cells.Item(109, 1).Value = 11
$ws.Cells.Item(109, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(109, 3).Value = "Bíobío"
$ws.Cells.Item(109, 4).Value = 44509
$ws.Cells.Item(109, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(109, 5).Value = 8
$ws.Cells.Item(109, 6).Value = "Fruta"
$ws.Cells.Item(109, 7).Value = 100108
$ws.Cells.Item(109, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(109, 9).Value = 100108005
$ws.Cells.Item(109, 10).Value = "Piña"
$ws.Cells.Item(109, 11).Value = "Caramelo"
$ws.Cells.Item(109, 12).Value = "Primera"
$ws.Cells.Item(109, 13).Value = 250
$ws.Cells.Item(109, 14).Value = 15000
$ws.Cells.Item(109, 15).Value = 16000
$ws.Cells.Item(109, 16).Value = 15520
$ws.Cells.Item(109, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(109, 18).Value = "Ecuador"
$ws.Cells.Item(109, 19).Value = 1109
$ws.Cells.Item(109, 20).Value = 14
